# Applies the "memory corrections" edit described in the commit:
#  - "Memory (MB)"!F2:F8 and F10:F25 become computed formulas (were hard-coded
#    literals) based on the "Datasets Attributes, Notes" sheet.
#  - A previously-skipped blank row 9 gets a styled (but empty) F9 cell.
#  - F20 changes from an integer number format to the same 0.00 format as the
#    rest of column F (picked up automatically because it now holds the same
#    formula as its neighbours).
#  - A handful of sheet selections / the active sheet change (pure view
#    state); the downstream "Memory - Input Ratio" sheet's cached formula
#    results shift accordingly once recalculated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Memory (MB)": replace the hard-coded F column values with formulas
# ---------------------------------------------------------------------
$wsMem = $wb.Worksheets.Item("Memory (MB)")

$memRows = 2,3,4,5,6,7,8,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25
foreach ($r in $memRows) {
    $formula = "=(1.5*'Datasets Attributes, Notes'!B$r*LOG('Datasets Attributes, Notes'!D$r,2) + 64 * 'Datasets Attributes, Notes'!D$r + 64)*0.000000125"
    $wsMem.Range("F$r").Formula = $formula
}

# Row 9 stays empty/unused, but gains a styled F9 cell matching the rest of
# the column (same 0.00 number format as F8/F10).
$wsMem.Range("F9").NumberFormat = "0.00"

# F20 previously used an integer style (s="4"); bring it in line with the
# rest of column F (s="5", "0.00") now that it also holds the shared formula.
$wsMem.Range("F20").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 2. Chart "Chart 5" (chart3.xml) plots 'Memory - Input Ratio'!F2:F8 (the
#    "subSeq" series) which is derived from the cells above - nudge Excel
#    to refresh its cached values now that the source data changed.
# ---------------------------------------------------------------------
$wsCharts = $wb.Worksheets.Item("Charts")
try {
    $chart5 = $wsCharts.ChartObjects("Chart 5").Chart
    $chart5.Refresh()
} catch {
}

$excel.CalculateFull()

# ---------------------------------------------------------------------
# 3. View-state: selections on a few sheets, and the active tab moves from
#    "Charts" to "Memory - Input Ratio".
# ---------------------------------------------------------------------

# "Memory (MB)": selection moves to F13 (sheet not left active).
$wsMem.Activate()
$wsMem.Range("F13").Select()

# "Datasets Attributes, Notes": row 17 gets selected (sheet not left active).
$wsNotes = $wb.Worksheets.Item("Datasets Attributes, Notes")
$wsNotes.Activate()
$wsNotes.Rows.Item(17).Select()

# "Memory - Input Ratio" becomes the active sheet/tab with F2 selected.
$wsRatio = $wb.Worksheets.Item("Memory - Input Ratio")
$wsRatio.Activate()
$wsRatio.Range("F2").Select()
